$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial values (column A) for the new trading-day rows being appended.
$dates = @(
    45838,
    45839,
    45840,
    45841,
    45845,
    45846,
    45847,
    45848,
    45849,
    45852,
    45853,
    45854,
    45855,
    45856,
    45859,
    45860,
    45861,
    45862,
    45863,
    45866,
    45867,
    45868,
    45869,
    45870,
    45873,
    45874,
    45875,
    45876,
    45877,
    45880,
    45881,
    45882,
    45883,
    45884,
    45887,
    45888,
    45889,
    45890,
    45891,
    45894,
    45895,
    45896,
    45897,
    45898,
    45902,
    45903,
    45904,
    45905,
    45908,
    45909,
    45910,
    45911,
    45912,
    45915
)

$startRow = 5018
$endRow = 5071
$count = $dates.Length

# Build a (count x 1) array for the bulk Value write.
$dateArr = New-Object 'object[,]' $count,1
for ($i = 0; $i -lt $count; $i++) {
    $dateArr[$i, 0] = $dates[$i]
}

$colA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$colB = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($endRow, 2))

$colA.Value = $dateArr
$colA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Sentiment column (B) stays at 0 for every new row, default (unstyled) cell format.
$colB.Value = 0

